# TrialsSetup workbook update: add a new trial row ("REMASTER (CLOU)") to the
# Power-Query-backed table on Sheet1, extending the table/range from
# A1:B15 to A1:B16, and keep the ExternalData_1 defined name in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow the query table by one row (this also extends the table ref/autoFilter
# and the sheet's used-range dimension to A1:B16).
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Populate the new row's "Trial Name" column; "Progress" (column B) is left
# blank, matching the source edit.
$ws.Range("A16").Value = "REMASTER (CLOU)"

# Keep the hidden ExternalData_1 defined name (used by the query table)
# pointing at the full, newly-extended range.
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$B`$16"
